$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; this shifts the existing rows 9..128 down to 10..129,
# matching the target diff (dimension grows from A1:R128 to A1:R129).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with its data. Columns A, B, C, E, F, G, H, I, N, Q, R
# keep the same values as the record that used to occupy row 9 (now shifted to row 10),
# while D, J, K, L, M, O, P take the new values specified by the edit.
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Macroferia Regional de Talca"
$ws.Range("C9").Value = "Maule"
$ws.Range("D9").Value = 44532
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 100112024
$ws.Range("G9").Value = "Choclo"
$ws.Range("H9").Value = "Choclero"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 30000
$ws.Range("K9").Value = 400
$ws.Range("L9").Value = 400
$ws.Range("M9").Value = 400
$ws.Range("N9").Value = "$/unidad"
$ws.Range("O9").Value = "Región de O'Higgins"
$ws.Range("P9").Value = 400
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = "Hortaliza"
